# Esercitazione2/Tabelle_Strategia2.xlsx
# "Completata Strategia 3 / Aggiungere seconda strategia per la somma"
#
# The only real content edit is the serial execution time for N=1000
# (cell C3 on "Foglio1"), which drives the first row of every
# "speed-up" / "efficiency" table (D11/E11, D19/E19, D28/E28) via the
# existing formulas (=C3/C11, =D11/2, etc.) - those recalc automatically.
# The new value is also displayed with one extra decimal digit
# (format "0.000" instead of "0.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$c3 = $ws.Range("C3")
$c3.Value = 1.4983
$c3.NumberFormat = "0.000"

# Reflect that the user ended up with C3 selected after editing it.
$c3.Select()
